# Greece Super League 1 - odds base update (06-04-2024 15:39)
# Swaps the home/away data pairs for a couple of already-played fixtures
# (rows 15/16 and 194/195), refreshes the odds for three upcoming
# fixtures (rows 199-201) with newer data, and drops four fixtures that
# are no longer part of the feed (old rows 202-205).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 (match id 13) ---------------------------------------------
$ws.Range("B15").Value  = 7100661
$ws.Range("F15").Value  = "Aris Salonika"
$ws.Range("G15").Value  = "Asteras Tripolis"
$ws.Range("H15").Value  = 3
$ws.Range("I15").Value  = 2
$ws.Range("K15").Value  = 1.8
$ws.Range("L15").Value  = 3.4
$ws.Range("M15").Value  = 4.75
$ws.Range("N15").Value  = 1.55
$ws.Range("O15").Value  = 3.8
$ws.Range("P15").Value  = 7
$ws.Range("Q15").Value  = -1
$ws.Range("R15").Value  = 2
$ws.Range("S15").Value  = 1.85
$ws.Range("T15").Value  = 2.25
$ws.Range("W15").Value  = 0.55
$ws.Range("Z15").Value  = 0
$ws.Range("AA15").Value = -0

# --- Row 16 (match id 14) ---------------------------------------------
$ws.Range("B16").Value  = 7100664
$ws.Range("F16").Value  = "Olympiakos"
$ws.Range("G16").Value  = "Lamia"
$ws.Range("H16").Value  = 4
$ws.Range("I16").Value  = 0
$ws.Range("K16").Value  = 1.25
$ws.Range("L16").Value  = 5.5
$ws.Range("M16").Value  = 12
$ws.Range("N16").Value  = 1.222
$ws.Range("O16").Value  = 6
$ws.Range("P16").Value  = 15
$ws.Range("Q16").Value  = -1.75
$ws.Range("R16").Value  = 1.875
$ws.Range("S16").Value  = 1.975
$ws.Range("T16").Value  = 2.75
$ws.Range("W16").Value  = 0.222
$ws.Range("Z16").Value  = 0.875
$ws.Range("AA16").Value = -1

# --- Row 194 (match id 192) --------------------------------------------
$ws.Range("B194").Value  = 7920471
$ws.Range("F194").Value  = "Aris Salonika"
$ws.Range("G194").Value  = "Lamia"
$ws.Range("H194").Value  = 3
$ws.Range("I194").Value  = 1
$ws.Range("K194").Value  = 1.571
$ws.Range("L194").Value  = 4
$ws.Range("M194").Value  = 6
$ws.Range("N194").Value  = 1.444
$ws.Range("O194").Value  = 4.5
$ws.Range("P194").Value  = 8.5
$ws.Range("Q194").Value  = -1.25
$ws.Range("R194").Value  = 1.925
$ws.Range("S194").Value  = 1.925
$ws.Range("T194").Value  = 2.75
$ws.Range("W194").Value  = 0.444
$ws.Range("Z194").Value  = 0.925
$ws.Range("AB194").Value = 1.025
$ws.Range("AC194").Value = -1

# --- Row 195 (match id 193) --------------------------------------------
$ws.Range("B195").Value  = 7920470
$ws.Range("F195").Value  = "AEK Athens"
$ws.Range("G195").Value  = "Olympiakos"
$ws.Range("H195").Value  = 1
$ws.Range("I195").Value  = 0
$ws.Range("K195").Value  = 1.909
$ws.Range("L195").Value  = 3.4
$ws.Range("M195").Value  = 4.2
$ws.Range("N195").Value  = 2.2
$ws.Range("O195").Value  = 3.2
$ws.Range("P195").Value  = 3.5
$ws.Range("Q195").Value  = -0.25
$ws.Range("R195").Value  = 1.85
$ws.Range("S195").Value  = 2
$ws.Range("T195").Value  = 2.5
$ws.Range("W195").Value  = 1.2
$ws.Range("Z195").Value  = 0.8500000000000001
$ws.Range("AB195").Value = -1
$ws.Range("AC195").Value = 0.825

# --- Row 199 (match id 197) refreshed with newer feed data ------------
$ws.Range("B199").Value = 7920478
$ws.Range("E199").Value = 45389.45833333334
$ws.Range("F199").Value = "Lamia"
$ws.Range("G199").Value = "Olympiakos"
$ws.Range("K199").Value = 15
$ws.Range("L199").Value = 6
$ws.Range("M199").Value = 1.2
$ws.Range("N199").Value = 10
$ws.Range("O199").Value = 5.75
$ws.Range("P199").Value = 1.285
$ws.Range("Q199").Value = 1.5
$ws.Range("R199").Value = 2.025
$ws.Range("S199").Value = 1.825
$ws.Range("T199").Value = 2.75
$ws.Range("U199").Value = 1.85
$ws.Range("V199").Value = 2

# --- Row 200 (match id 198) refreshed with newer feed data ------------
$ws.Range("B200").Value = 7920477
$ws.Range("E200").Value = 45389.5625
$ws.Range("F200").Value = "Aris Salonika"
$ws.Range("G200").Value = "Panathinaikos"
$ws.Range("K200").Value = 3.75
$ws.Range("L200").Value = 3.4
$ws.Range("M200").Value = 2
$ws.Range("N200").Value = 4
$ws.Range("O200").Value = 3.5
$ws.Range("P200").Value = 1.909
$ws.Range("Q200").Value = 0.5
$ws.Range("R200").Value = 1.925
$ws.Range("S200").Value = 1.925
$ws.Range("U200").Value = 1.825
$ws.Range("V200").Value = 2.025

# --- Row 201 (match id 199) refreshed with newer feed data ------------
$ws.Range("B201").Value = 7920476
$ws.Range("E201").Value = 45389.60416666666
$ws.Range("F201").Value = "AEK Athens"
$ws.Range("G201").Value = "PAOK Salonika"
$ws.Range("K201").Value = 1.909
$ws.Range("L201").Value = 3.4
$ws.Range("M201").Value = 4
$ws.Range("N201").Value = 1.75
$ws.Range("O201").Value = 3.5
$ws.Range("P201").Value = 5
$ws.Range("Q201").Value = -0.75
$ws.Range("R201").Value = 1.975
$ws.Range("S201").Value = 1.875
$ws.Range("U201").Value = 1.975
$ws.Range("V201").Value = 1.875

# --- Drop the four fixtures that fell off the feed (old rows 202-205) --
$ws.Rows("202:205").Delete()
